# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 6.201049113329182

# Row 3
$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 3.781711156805759

# Row 4
$ws.Range("B4").Value = 0.04763786555579896
$ws.Range("C4").Value = 114.8270160096505
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 124.3403759977393

# Row 5
$ws.Range("B5").Value = 1.459612070389937
$ws.Range("C5").Value = 10.29869402782916
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 21.22402822075207

# Row 6
$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 3099.503889238888
$ws.Range("D6").Value = 0.8054896365839992
$ws.Range("E6").Value = 8.660232485948974
$ws.Range("G6").Value = 3112.200597044728

# Row 7
$ws.Range("B7").Value = 1.459612070389937
$ws.Range("C7").Value = 10.29869402782916
$ws.Range("D7").Value = 0.8054896365839992
$ws.Range("E7").Value = 645.3272768299601
$ws.Range("G7").Value = 657.8910725647631
